# "Generate Report for Handoff"
#
# Updates the localization-status report to reflect that the content is
# now ready for handoff (rather than still in translation) and refreshes
# the generation timestamps, then widens the "Latest Handoff/Handback
# Datetime"-ish columns that now need to fit the longer status text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" -----------------
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Refreshed "generate date" / handoff datetime timestamps --------------
$overview.Range("G2").Value = "2016-09-02 12:45:31"
$zhcn.Range("H2").Value     = "2016-09-02 12:45:27"
$dede.Range("H2").Value     = "2016-09-02 12:45:31"

# --- Widen the status columns to fit "Ready for handoff" -------------------
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth     = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth     = 16.3333333333333
